$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for Cilantro at the top of the
# series (Macroferia Regional de Talca), pushing the existing rows
# (16-47) down by one (to 17-48) and growing the used range to R48.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new record's values.
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44791
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112040
$ws.Range("G16").Value = "Cilantro"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 150
$ws.Range("K16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = 10000
$ws.Range("N16").Value = "$/caja 36 atados"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 278
$ws.Range("Q16").Value = 36
$ws.Range("R16").Value = "Hortaliza"
